# Update the "Förändrad" date column (C2:C15) from 2023-09-19 (45188)
# to 2023-09-20 (45189), matching the automatic update reflected in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
